$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.4375,
    0.421875,
    0.390625,
    0.296875,
    0.3125,
    0.3125,
    0.34375,
    0.328125,
    0.3125,
    0.296875,
    0.265625,
    0.3125,
    0.28125,
    0.234375,
    0.28125,
    0.265625,
    0.390625,
    0.3125,
    0.234375,
    0.25,
    0.3125,
    0.328125,
    0.28125,
    0.265625,
    0.3125,
    0.25,
    0.265625,
    0.28125,
    0.25,
    0.28125,
    0.25,
    0.25,
    0.25,
    0.25,
    0.234375,
    0.234375,
    0.25,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.265625,
    0.109375,
    0.15625,
    0.234375,
    0.140625,
    0.25,
    0.265625,
    0.203125,
    0.203125,
    0.171875,
    0.171875,
    0.234375,
    0.171875,
    0.171875,
    0.15625,
    0.28125,
    0.2295081967213115
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("A2:B118").Select()
